# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 537
$ws1.Range("F5").Value = 520
$ws1.Range("F7").Value = 2626
$ws1.Range("F8").Value = 448
$ws1.Range("F9").Value = 7206
$ws1.Range("F10").Value = 191
$ws1.Range("F11").Value = 449
$ws1.Range("F13").Value = 163
$ws1.Range("F14").Value = 38

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 537
$ws4.Range("F5").Value = 520
$ws4.Range("F9").Value = 2626
$ws4.Range("F10").Value = 448
$ws4.Range("F11").Value = 7206
$ws4.Range("F12").Value = 191
$ws4.Range("F13").Value = 449
$ws4.Range("F17").Value = 163
$ws4.Range("F18").Value = 38
